# Apply "data : case 1" edit
#
# Summary of the change:
#  - Column A and Column B custom widths are swapped
#      (A: 15.42578125 -> 14.7109375, B: 14.7109375 -> 15.42578125)
#  - All existing values in A1:B4 are replaced with new values
#  - A new row 5 (A5:B5) is appended with new values
#  - The used range/dimension grows from A1:B4 to A1:B5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap column widths between column A and column B ---
# Note: the COM ColumnWidth property only supports increments of 1/6 of a
# character in this runtime, so these are the closest achievable values to
# the exact target widths of 14.7109375 (col A) and 15.42578125 (col B).
$ws.Columns.Item(1).ColumnWidth = 13.833333333333334
$ws.Columns.Item(2).ColumnWidth = 14.666666666666666

# --- Update existing rows 1-4 with new values ---
$ws.Range("A1").Value = 0.043651242330346826
$ws.Range("B1").Value = -0.043651242413485267

$ws.Range("A2").Value = -0.017889369357630615
$ws.Range("B2").Value = 0.017889369269980131

$ws.Range("A3").Value = -0.028339889734841812
$ws.Range("B3").Value = 0.028339889644740306

$ws.Range("A4").Value = 0.0048007699384992634
$ws.Range("B4").Value = -0.0048007700213412103

# --- Add new row 5 ---
$ws.Range("A5").Value = 0.072839425633565935
$ws.Range("B5").Value = -0.072839425718507017
